$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.424577333333334
$ws.Range("H2").Value = 4.273732000000001
$ws.Range("I2").Value = 0.009249507402003717
$ws.Range("J2").Value = 0.009249507402003719
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.02725333333333333
$ws.Range("N2").Value = 0.08176
$ws.Range("O2").Value = 0.0007089206372884383
$ws.Range("P2").Value = 0.0007089206372884382
$ws.Range("Q2").Value = 0.03882448092444445
$ws.Range("R2").Value = 0.3494203283200001
$ws.Range("S2").Value = 0.000006557166682032603
$ws.Range("T2").Value = 0.000006557166682032603

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.424577333333334
$ws.Range("H3").Value = 4.273732000000001
$ws.Range("I3").Value = 0.009249507402003717
$ws.Range("J3").Value = 0.009249507402003719
$ws.Range("M3").Value = 38.416166
$ws.Range("N3").Value = 115.248498
$ws.Range("O3").Value = 0.9992910793627116
$ws.Range("P3").Value = 0.9992910793627116
$ws.Range("Q3").Value = 54.72679931717067
$ws.Range("R3").Value = 492.5411938545361
$ws.Range("S3").Value = 0.009242950235321686
$ws.Range("T3").Value = 0.009242950235321687

$ws.Range("I4").Value = 0.6303829267608616
$ws.Range("J4").Value = 0.6303829267608616
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.02725333333333333
$ws.Range("N4").Value = 0.08176
$ws.Range("O4").Value = 0.0007089206372884383
$ws.Range("P4").Value = 0.0007089206372884382
$ws.Range("Q4").Value = 2.646010090204444
$ws.Range("R4").Value = 23.81409081184
$ws.Range("S4").Value = 0.0004468914661750609
$ws.Range("T4").Value = 0.0004468914661750609

$ws.Range("I5").Value = 0.6303829267608616
$ws.Range("J5").Value = 0.6303829267608616
$ws.Range("M5").Value = 38.416166
$ws.Range("N5").Value = 115.248498
$ws.Range("O5").Value = 0.9992910793627116
$ws.Range("P5").Value = 0.9992910793627116
$ws.Range("Q5").Value = 3729.802942623614
$ws.Range("R5").Value = 33568.22648361253
$ws.Range("S5").Value = 0.6299360352946866
$ws.Range("T5").Value = 0.6299360352946866

$ws.Range("G6").Value = 31.31606233333333
$ws.Range("H6").Value = 93.94818699999999
$ws.Range("I6").Value = 0.2033291865426586
$ws.Range("J6").Value = 0.2033291865426586
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.02725333333333333
$ws.Range("N6").Value = 0.08176
$ws.Range("O6").Value = 0.0007089206372884383
$ws.Range("P6").Value = 0.0007089206372884382
$ws.Range("Q6").Value = 0.8534670854577777
$ws.Range("R6").Value = 7.681203769119999
$ws.Range("S6").Value = 0.0001441442565031613
$ws.Range("T6").Value = 0.0001441442565031612

$ws.Range("G7").Value = 31.31606233333333
$ws.Range("H7").Value = 93.94818699999999
$ws.Range("I7").Value = 0.2033291865426586
$ws.Range("J7").Value = 0.2033291865426586
$ws.Range("M7").Value = 38.416166
$ws.Range("N7").Value = 115.248498
$ws.Range("O7").Value = 0.9992910793627116
$ws.Range("P7").Value = 0.9992910793627116
$ws.Range("Q7").Value = 1203.04304906368
$ws.Range("R7").Value = 10827.38744157312
$ws.Range("S7").Value = 0.2031850422861554
$ws.Range("T7").Value = 0.2031850422861554

$ws.Range("G8").Value = 24.18651133333333
$ws.Range("H8").Value = 72.559534
$ws.Range("I8").Value = 0.1570383792944762
$ws.Range("J8").Value = 0.1570383792944762
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.02725333333333333
$ws.Range("N8").Value = 0.08176
$ws.Range("O8").Value = 0.0007089206372884383
$ws.Range("P8").Value = 0.0007089206372884382
$ws.Range("Q8").Value = 0.6591630555377778
$ws.Range("R8").Value = 5.93246749984
$ws.Range("S8").Value = 0.0001113277479281835
$ws.Range("T8").Value = 0.0001113277479281836

$ws.Range("G9").Value = 24.18651133333333
$ws.Range("H9").Value = 72.559534
$ws.Range("I9").Value = 0.1570383792944762
$ws.Range("J9").Value = 0.1570383792944762
$ws.Range("M9").Value = 38.416166
$ws.Range("N9").Value = 115.248498
$ws.Range("O9").Value = 0.9992910793627116
$ws.Range("P9").Value = 0.9992910793627116
$ws.Range("Q9").Value = 929.1530343422145
$ws.Range("R9").Value = 8362.377309079931
$ws.Range("S9").Value = 0.156927051546548
$ws.Range("T9").Value = 0.156927051546548
